# Fruta / hortaliza, semanal
# Insert a new weekly observation as row 93 (shifting the existing rows
# 93..190 down to 94..191), matching the new canonical dimension A1:T191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 93:190 down by inserting a new blank row at 93.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(93, 1).Value = 9
$ws.Cells.Item(93, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(93, 3).Value = "Metropolitana"
$ws.Cells.Item(93, 4).Value = 45280
$ws.Cells.Item(93, 5).Value = 13
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100101
$ws.Cells.Item(93, 8).Value = "Berries"
$ws.Cells.Item(93, 9).Value = 100101004
$ws.Cells.Item(93, 10).Value = "Frambuesa"
$ws.Cells.Item(93, 11).Value = "Sin especificar"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 600
$ws.Cells.Item(93, 14).Value = 6000
$ws.Cells.Item(93, 15).Value = 6000
$ws.Cells.Item(93, 16).Value = 6000
$ws.Cells.Item(93, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(93, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(93, 19).Value = 3000
$ws.Cells.Item(93, 20).Value = 2
